$d = $word.ActiveDocument

$replacements = @(
    @("87÷7=", "19÷5="),
    @("30÷5=", "19÷2="),
    @("87÷9=", "95÷5="),
    @("18÷2=", "96÷2="),
    @("17÷4=", "68÷7="),
    @("51÷9=", "76÷2="),
    @("76÷5=", "23÷3="),
    @("52÷7=", "23÷2="),
    @("25÷5=", "76÷3="),
    @("21÷8=", "83÷4="),
    @("68÷3=", "48÷3="),
    @("86÷4=", "20÷9="),
    @("60÷9=", "99÷3="),
    @("72÷6=", "35÷7="),
    @("33÷2=", "49÷7="),
    @("94÷7=", "13÷6="),
    @("32÷3=", "19÷9="),
    @("31÷6=", "84÷3="),
    @("74÷8=", "91÷5="),
    @("27÷8=", "66÷7="),
    @("19÷7=", "51÷5="),
    @("60÷2=", "23÷6="),
    @("25÷2=", "47÷3="),
    @("14÷9=", "64÷8="),
    @("13÷2=", "80÷8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
